# Update "Generate Report for Handback" timestamps.
#
# The workbook has three sheets:
#   - Overview : G2 = "Latest HO Xliff Generate Date" value
#   - zh-cn    : H2 = "Correspond Handoff Datetime", K2 = "Correspond Handback DateTime"
#   - de-de    : H2 = "Correspond Handoff Datetime", K2 = "Correspond Handback DateTime"
#
# Overview!G2 and de-de!H2 originally shared the same text value
# ("2016-09-06 03:13:00"), so both are refreshed to the same new value.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Keep these as plain text so Excel doesn't reinterpret them as date serials.
$wsOverview.Range("G2").Value = "2016-09-06 03:13:47"

$wsZhCn.Range("H2").Value = "2016-09-06 03:13:43"
$wsZhCn.Range("K2").Value = "2016-09-06 03:14:01"

$wsDeDe.Range("H2").Value = "2016-09-06 03:13:47"
$wsDeDe.Range("K2").Value = "2016-09-06 03:14:15"
